# Apply the changes described in the commit:
#   "Add further tests for NetworkCost class. Reorganize modules.
#    Rework classes to proximately transfer responsibilities to
#    FreightNetwork class."
#
# Concretely, in the "railway_parameters" workbook this means adding two
# new parameter rows to the "mob" sheet (net_to_gross_factor and
# main_min_density), and switching the active sheet/selection from "inf"
# back to "mob".

$wb = $excel.ActiveWorkbook

$mob = $wb.Worksheets.Item("mob")
$inf = $wb.Worksheets.Item("inf")

# --- Add the two new rows to the "mob" sheet ---------------------------

# Row 23: net_to_gross_factor
$mob.Range("A23").Value = "net_to_gross_factor"
$mob.Range("B23").Value = 1.67
$mob.Range("C23").Value = "Factor to convert (aproximately) net tons to gross tons, based on a full train operation (coef)."

# Row 24: main_min_density (rendered with the same thousands-style number
# format used by the other "big number" cells in the column, e.g. B8/B9)
$mob.Range("A24").Value = "main_min_density"
$mob.Range("B24").Value = 760000
$mob.Range("B24").NumberFormat = "#,##0"
$mob.Range("C24").Value = "Minimum density to consider a link as being a main track (ton-km/ton = ton)."

# --- Move the active sheet/selection from "inf" back to "mob" ----------

# Previously "inf" was the active/selected tab with A12 selected; now
# "mob" becomes active again, scrolled down and with the newly added
# row selected.
$mob.Activate()
$mob.Range("A24:C24").Select()

# Leave "inf"'s own selection (A12) untouched - only the tab-selection
# flag moves because "mob" is activated above.
